$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Single Page Architecture:"
# and the (empty, bulleted) paragraph immediately before it - that is
# the third bullet under "Links:" which is about to receive the new
# YouTube hyperlink.
$headingPara = $null
$linkPara = $null
$prevPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Single Page Architecture:") {
        $headingPara = $p
        $linkPara = $prevPara
    }
    $prevPara = $p
}

# --- 1. Turn that empty bullet into a hyperlink pointing at the new
#        YouTube video -------------------------------------------------
$d.Hyperlinks.Add($linkPara.Range, "https://www.youtube.com/watch?v=RWXKysImabs")

# --- 2. The "Single Page Architecture:" paragraph is emptied out and
#        becomes another bullet in the same list -----------------------
$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr></w:p></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
$null = $headingPara.Range.InsertXML($newParaXml)
